$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.434.10"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.641.49"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.50"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.73"
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("D9").Value = "2.643.22"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.49"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "3.114.10"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "63.403.41"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000145"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "2.623.83"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.34"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.01"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.69"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.56"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.69"
$ws.Range("E25").Value = "  +6.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.55"
$ws.Range("E26").Value = "  +6.33%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "555.52"
$ws.Range("E28").Value = "  +17.62%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.43"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.82"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  +12.33%  "
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").Value = "0.0₃0802"
$ws.Range("E34").Value = "  -2.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "174.99"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("E36").Value = "  +7.35%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.401"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.03"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "169.97"
$ws.Range("E42").Value = "  +7.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.27"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.71"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.20"
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0549"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0957"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0237"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.68"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.35"
$ws.Range("E51").Value = "  -0.61%  "
